# MVP V2 Scrape Inventory and Apply Purchase Amount at checkout
#
# The scraped inventory table on sheet "0005" (A1:B6, "Item"/"Quantity")
# is re-sorted alphabetically (ascending) by the Item column - i.e. the
# same effect as running Data > Sort on A2:B6 keyed on column A, with no
# header row included in the sorted range. This is what Excel itself does
# (and records) when a user applies a standard A-Z sort from the ribbon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Sort.SortFields.Clear()
# xlSortOnValues=0, xlAscending=1, xlSortNormal=0
$ws.Sort.SortFields.Add2($ws.Range("A2:A6"), 0, 1, $null, 0)
$ws.Sort.SetRange($ws.Range("A2:B6"))
$ws.Sort.Header = 2          # xlNo - the range being sorted excludes the header row
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1     # xlSortColumns - sort top-to-bottom by column
$ws.Sort.SortMethod = 1      # xlPinYin
$ws.Sort.Apply()
